$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "249.58"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.78"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.439"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05615"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.426"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.371"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8123"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8964"

# Row 10
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01149"
$ws.Range("E10").Value = "9OneONEBestin24h"

# Row 11
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1430"
$ws.Range("E11").Value = "10WazirXWRX"

# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07532"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

# Row 13
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03094"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

# Row 14
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03092"
$ws.Range("E14").Value = "13BitrueCoinBTR"

# Row 15
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09333"
$ws.Range("E15").Value = "14BitMartTokenBMX"

# Row 16
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.561"
$ws.Range("E16").Value = "15MCDexMCB"

# Row 17
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001592"
$ws.Range("E17").Value = "16BitForexTokenBF"

# Row 18
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04755"
$ws.Range("E18").Value = "17CoinExTokenCET"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006417"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004988"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.188"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3304"

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0003010"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04020"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006820"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1066"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002741"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007781"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005580"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5007"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2396"
$ws.Range("E48").Value = "47BOLOBOLO"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002103"
